# "little bug fix Optimization Done!!"
#
# Row 15 (10:06 -> 10:14, 8 min) gets corrected/split into three rows:
#   Row 15: 10:06 -> 10:18 (12 min)   [end time bug-fixed from 10:14 to 10:18]
#   Row 16: 10:19 -> 10:27 (8 min)    [new row - what used to be row 15's slot]
#   Row 17: 09:30 -> 10:10 (40 min)   [new row, different category]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 15: end time 10:14 -> 10:18, duration 8 -> 12
$ws.Cells.Item(15, 1).Value = 19
$ws.Cells.Item(15, 2).Value = "10:06"
$ws.Cells.Item(15, 3).Value = "10:18"
$ws.Cells.Item(15, 4).Value = 12
$ws.Cells.Item(15, 5).Value = 3

# Insert new row 16
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "10:19"
$ws.Cells.Item(16, 3).Value = "10:27"
$ws.Cells.Item(16, 4).Value = 8
$ws.Cells.Item(16, 5).Value = 3

# Insert new row 17
$ws.Cells.Item(17, 1).Value = 20
$ws.Cells.Item(17, 2).Value = "09:30"
$ws.Cells.Item(17, 3).Value = "10:10"
$ws.Cells.Item(17, 4).Value = 40
$ws.Cells.Item(17, 5).Value = 1
